$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generales")
$c = $ws.Range("J1").Comment
Write-Host "Parent type test"
try {
  $c.Parent = $ws.Range("L1")
  Write-Host "set parent ok"
} catch { Write-Host "cant set parent:" $_.Exception.Message }
Write-Host $ws.Range("L1").Comment
